$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for account 004491730 (Denise) first so that the earlier
# row number (3) for account 004450724 (Assako) stays valid.
$ws.Rows.Item(227).Delete()
$ws.Rows.Item(3).Delete()
